$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users")

$newUserIds = @("9384702a", "06e70e3a", "c3b16923", "3cba4514", "2202a7bd", "c6240f5e", "cd247874", "165f6259", "75e4d446", "03bcf857", "63298ed8", "a1fb7d9b", "5dccde57", "4525f6ee", "da83814f", "80b4c93b", "53a1f3d2", "2cde6e96", "cf308c06", "0b54e04a", "6778b084", "e4abec67", "c24ad56c", "50ce42f7", "7da8bde3", "a58e73db", "70430bdd", "f9ea08e2", "03d3834b", "f7d7c5b7", "62991d2e", "eae72c4d", "9d6ce120", "534e98e4", "fe5c8ae6", "3f6c8015", "c89879cc", "91254d01", "fc890919", "510cb503", "d98fdbf7", "a3549972", "0795cd3e", "e8e60d2e", "a26dd3b6", "57b476a7", "359cf14d", "bf3a8843", "89c2d2f6", "2e492def")

for ($i = 0; $i -lt $newUserIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newUserIds[$i]
}

Write-Host "Updated $($newUserIds.Length) User ID values in column A"
